# Commit: "remove description from test xlsx file."
#
# The constraint table lives on worksheet "B" and has columns:
#   A=constraint, B=description, C=dir, D=rhs
# The "description" column (B) is removed, shifting dir/rhs left
# (C->B, D->C). The author also left the workbook with sheet "B" as the
# active/selected sheet, with column B (now "dir") selected.

$wb = $excel.ActiveWorkbook
$wsB = $wb.Worksheets.Item("B")

# Drop the whole "description" column; C (dir) and D (rhs) shift left.
$wsB.Columns("B").Delete()

# Make "B" the active sheet and leave column B (the new "dir" column)
# selected, matching the saved selection state.
$wsB.Activate()
$null = $wsB.Columns("B").Select()
